# Generate Report for Handoff
# The localization run produced a new handoff package for "b.md" (in both the
# zh-cn and de-de target sheets) and flagged that its previously recorded
# handback file is now stale. Update the Overview summary sheet and both
# locale detail sheets accordingly, and widen the "Error Detail" column so
# the new warning text is readable.

$wb = $excel.ActiveWorkbook

$newHandoffDate = "2016-08-21 04:44:41"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87c4c5cf39de093bbcba597b26e62bf011a5a608/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca24de754c41a0e45c841f01415c6aa8637eeb6e/e2e/b.md."

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is "b.md". Both locale status columns move to
# "Ready for handoff" and the latest-generate-date column is refreshed.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $newHandoffDate

# ---------------------------------------------------------------------------
# zh-cn detail sheet: row 3 ("b.md") gets the new handoff file/date, its
# "Content Duplicate" flag clears, status flips, and the error detail is
# populated. Column P ("Error Detail") is widened to fit the message.
#
# Note: writing the literal text "False" straight into a Range.Value makes
# Excel coerce it to a real boolean cell (t="b"), but the source file always
# stores these flags as plain text (t="s"). Prefixing with an apostrophe
# forces text entry like Excel's own "quote prefix", then resetting the
# style back to Normal drops the quote-prefix formatting it leaves behind.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-21 04:44:37"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de detail sheet: same shape of change as zh-cn.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $newHandoffDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
